$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ur = $ws.UsedRange
Write-Host "UsedRange Address:" $ur.Address()
Write-Host "Rows:" $ur.Rows.Count
Write-Host "Cols:" $ur.Columns.Count
